# Commit: "Added test to write data to excel and json"
# Appends three new rows of email/password test data below the existing
# table (rows 1-3 stay untouched), extending the used range to A1:B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "test@example.com"
$ws.Range("B4").Value = "password123"

$ws.Range("A5").Value = "test@example.com"
$ws.Range("B5").Value = "password123"

$ws.Range("A6").Value = "tester@ample.com"
$ws.Range("B6").Value = "password754"
